# Auto-generated Excel COM-interop script applying the diff
$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 19:35:31"
$ws1.Range("A3").Value = "Total filas: 116"

$ws1.Cells.Item(19,1).Value = "16:50:41"
$ws1.Cells.Item(19,2).Value = "17:17"
$ws1.Cells.Item(19,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(19,4).Value = 27
$ws1.Cells.Item(19,5).Value = "LP1912"

$ws1.Cells.Item(20,1).Value = "16:46:42"
$ws1.Cells.Item(20,2).Value = "17:17"
$ws1.Cells.Item(20,3).Value = "17_ROMERO"
$ws1.Cells.Item(20,4).Value = 31
$ws1.Cells.Item(20,5).Value = "LP1912"

$ws1.Cells.Item(48,1).Value = "18:10:41"
$ws1.Cells.Item(48,2).Value = "18:11"
$ws1.Cells.Item(48,3).Value = "10_OLMOS"
$ws1.Cells.Item(48,4).Value = 1
$ws1.Cells.Item(48,5).Value = "LP1912"

$ws1.Cells.Item(49,1).Value = "18:10:41"
$ws1.Cells.Item(49,2).Value = "18:11"
$ws1.Cells.Item(49,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(49,4).Value = 1
$ws1.Cells.Item(49,5).Value = "LP1912"

$ws1.Cells.Item(76,1).Value = "18:44:34"
$ws1.Cells.Item(76,2).Value = "19:17"
$ws1.Cells.Item(76,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(76,4).Value = 33
$ws1.Cells.Item(76,5).Value = "LP1912"

$ws1.Cells.Item(77,1).Value = "17:47:22"
$ws1.Cells.Item(77,2).Value = "19:17"
$ws1.Cells.Item(77,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(77,4).Value = 90
$ws1.Cells.Item(77,5).Value = "LP1912"

$ws1.Cells.Item(95,1).Value = "19:35:31"
$ws1.Cells.Item(95,2).Value = "19:41"
$ws1.Cells.Item(95,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(95,4).Value = 6
$ws1.Cells.Item(95,5).Value = "LP1912"

$ws1.Cells.Item(96,1).Value = "19:11:59"
$ws1.Cells.Item(96,2).Value = "19:43"
$ws1.Cells.Item(96,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(96,4).Value = 32
$ws1.Cells.Item(96,5).Value = "LP1912"

$ws1.Cells.Item(97,1).Value = "17:47:22"
$ws1.Cells.Item(97,2).Value = "19:44"
$ws1.Cells.Item(97,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(97,4).Value = 117
$ws1.Cells.Item(97,5).Value = "LP1912"

$ws1.Cells.Item(98,1).Value = "18:31:18"
$ws1.Cells.Item(98,2).Value = "19:46"
$ws1.Cells.Item(98,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(98,4).Value = 75
$ws1.Cells.Item(98,5).Value = "LP1912"

$ws1.Cells.Item(99,1).Value = "19:11:59"
$ws1.Cells.Item(99,2).Value = "19:50"
$ws1.Cells.Item(99,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(99,4).Value = 39
$ws1.Cells.Item(99,5).Value = "LP1912"

$ws1.Cells.Item(100,1).Value = "17:54:43"
$ws1.Cells.Item(100,2).Value = "19:51"
$ws1.Cells.Item(100,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(100,4).Value = 117
$ws1.Cells.Item(100,5).Value = "LP1912"

$ws1.Cells.Item(101,1).Value = "18:10:41"
$ws1.Cells.Item(101,2).Value = "19:58"
$ws1.Cells.Item(101,3).Value = "14X44_ABASTO"
$ws1.Cells.Item(101,4).Value = 108
$ws1.Cells.Item(101,5).Value = "LP1912"

$ws1.Cells.Item(102,1).Value = "18:31:18"
$ws1.Cells.Item(102,2).Value = "19:59"
$ws1.Cells.Item(102,3).Value = "14X44_ABASTO"
$ws1.Cells.Item(102,4).Value = 88
$ws1.Cells.Item(102,5).Value = "LP1912"

$ws1.Cells.Item(103,1).Value = "18:10:41"
$ws1.Cells.Item(103,2).Value = "20:00"
$ws1.Cells.Item(103,3).Value = "215C_EL PATO"
$ws1.Cells.Item(103,4).Value = 110
$ws1.Cells.Item(103,5).Value = "LP1912"

$ws1.Cells.Item(104,1).Value = "18:31:18"
$ws1.Cells.Item(104,2).Value = "20:01"
$ws1.Cells.Item(104,3).Value = "215C_EL PATO"
$ws1.Cells.Item(104,4).Value = 90
$ws1.Cells.Item(104,5).Value = "LP1912"

$ws1.Cells.Item(105,1).Value = "19:11:59"
$ws1.Cells.Item(105,2).Value = "20:04"
$ws1.Cells.Item(105,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(105,4).Value = 53
$ws1.Cells.Item(105,5).Value = "LP1912"

$ws1.Cells.Item(106,1).Value = "19:35:31"
$ws1.Cells.Item(106,2).Value = "20:10"
$ws1.Cells.Item(106,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(106,4).Value = 35
$ws1.Cells.Item(106,5).Value = "LP1912"

$ws1.Cells.Item(107,1).Value = "19:11:59"
$ws1.Cells.Item(107,2).Value = "20:13"
$ws1.Cells.Item(107,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(107,4).Value = 62
$ws1.Cells.Item(107,5).Value = "LP1912"

$ws1.Cells.Item(108,1).Value = "18:31:18"
$ws1.Cells.Item(108,2).Value = "20:14"
$ws1.Cells.Item(108,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(108,4).Value = 103
$ws1.Cells.Item(108,5).Value = "LP1912"

$ws1.Cells.Item(109,1).Value = "19:11:59"
$ws1.Cells.Item(109,2).Value = "20:25"
$ws1.Cells.Item(109,3).Value = "15_ABASTO"
$ws1.Cells.Item(109,4).Value = 74
$ws1.Cells.Item(109,5).Value = "LP1912"

$ws1.Cells.Item(110,1).Value = "18:31:18"
$ws1.Cells.Item(110,2).Value = "20:26"
$ws1.Cells.Item(110,3).Value = "15_ABASTO"
$ws1.Cells.Item(110,4).Value = 115
$ws1.Cells.Item(110,5).Value = "LP1912"

$ws1.Cells.Item(111,1).Value = "18:44:34"
$ws1.Cells.Item(111,2).Value = "20:28"
$ws1.Cells.Item(111,3).Value = "10_OLMOS"
$ws1.Cells.Item(111,4).Value = 104
$ws1.Cells.Item(111,5).Value = "LP1912"

$ws1.Cells.Item(112,1).Value = "18:31:18"
$ws1.Cells.Item(112,2).Value = "20:29"
$ws1.Cells.Item(112,3).Value = "10_OLMOS"
$ws1.Cells.Item(112,4).Value = 118
$ws1.Cells.Item(112,5).Value = "LP1912"

$ws1.Cells.Item(113,1).Value = "19:11:59"
$ws1.Cells.Item(113,2).Value = "20:43"
$ws1.Cells.Item(113,3).Value = "215B_EL PATO"
$ws1.Cells.Item(113,4).Value = 92
$ws1.Cells.Item(113,5).Value = "LP1912"

$ws1.Cells.Item(114,1).Value = "19:11:59"
$ws1.Cells.Item(114,2).Value = "20:44"
$ws1.Cells.Item(114,3).Value = "17X38_ROMERO"
$ws1.Cells.Item(114,4).Value = 93
$ws1.Cells.Item(114,5).Value = "LP1912"

$ws1.Cells.Item(115,1).Value = "18:52:04"
$ws1.Cells.Item(115,2).Value = "20:44"
$ws1.Cells.Item(115,3).Value = "215B_EL PATO"
$ws1.Cells.Item(115,4).Value = 112
$ws1.Cells.Item(115,5).Value = "LP1912"

$ws1.Cells.Item(116,1).Value = "18:52:04"
$ws1.Cells.Item(116,2).Value = "20:45"
$ws1.Cells.Item(116,3).Value = "17X38_ROMERO"
$ws1.Cells.Item(116,4).Value = 113
$ws1.Cells.Item(116,5).Value = "LP1912"

$ws1.Cells.Item(117,1).Value = "19:35:31"
$ws1.Cells.Item(117,2).Value = "20:52"
$ws1.Cells.Item(117,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(117,4).Value = 77
$ws1.Cells.Item(117,5).Value = "LP1912"

$ws1.Cells.Item(118,1).Value = "19:11:59"
$ws1.Cells.Item(118,2).Value = "21:01"
$ws1.Cells.Item(118,3).Value = "215A_EL PATO"
$ws1.Cells.Item(118,4).Value = 110
$ws1.Cells.Item(118,5).Value = "LP1912"

$ws1.Cells.Item(119,1).Value = "19:11:59"
$ws1.Cells.Item(119,2).Value = "21:02"
$ws1.Cells.Item(119,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(119,4).Value = 111
$ws1.Cells.Item(119,5).Value = "LP1912"

$ws1.Cells.Item(120,1).Value = "19:35:31"
$ws1.Cells.Item(120,2).Value = "21:10"
$ws1.Cells.Item(120,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(120,4).Value = 95
$ws1.Cells.Item(120,5).Value = "LP1912"

$ws1.Cells.Item(121,1).Value = "19:35:31"
$ws1.Cells.Item(121,2).Value = "21:23"
$ws1.Cells.Item(121,3).Value = "10_OLMOS"
$ws1.Cells.Item(121,4).Value = 108
$ws1.Cells.Item(121,5).Value = "LP1912"

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 19:35:31"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 19:35:31"
$ws3.Range("A3").Value = "Total filas: 13"
$ws3.Cells.Item(18,1).Value = "19:35:31"
$ws3.Cells.Item(18,2).Value = "21:34"
$ws3.Cells.Item(18,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(18,4).Value = 119
$ws3.Cells.Item(18,5).Value = "L6203"
